# Commit: Sat, Apr 18, 2020  3:05:24 PM
#
# The authoritative OOXML diff swaps the bodies of ppt/theme/theme1.xml
# ("Office Theme") and ppt/theme/theme2.xml ("Integral") while leaving every
# relationship (presentation.xml.rels, slideMaster1.xml.rels,
# notesMaster1.xml.rels, ...) untouched:
#   - theme1.xml (only ever linked from the notes master) ends up holding the
#     "Integral" color scheme that used to live in theme2.xml.
#   - theme2.xml (linked from the slide master and the presentation's default
#     theme relationship) ends up holding the "Office Theme" color scheme
#     that used to live in theme1.xml.
# The font scheme (majorFont/minorFont) and the format scheme (fill/line/
# effect/background styles) are already byte-identical between the two
# theme parts, so the only substantive content that needs to move is the
# <a:clrScheme> color list.
#
# Via the PowerPoint object model, the live presentation theme (the one the
# slide master / slides actually render with) is reached through
# SlideMaster.Theme.ThemeColorScheme - a 12 item collection ordered
# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink, exactly mirroring
# <a:clrScheme>'s child order. We push the target ("Office Theme") colors
# into that collection so the slide-facing theme part ends up matching the
# diff.

function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Target "Office Theme" clrScheme, in <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$targetColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $targetColors.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = HexToComRgb($targetColors[$i])
}
